$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

# C1: keep only top/bottom border (drop left & right) -> borderId 4
$ws1.Range("C1").Borders.Item(7).LineStyle = 0
$ws1.Range("C1").Borders.Item(10).LineStyle = 0

# D1: keep top/right/bottom border (drop left) -> borderId 5
$ws1.Range("D1").Borders.Item(7).LineStyle = 0

# C2: rename "fedcore" column header to "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

# C1: keep only top/bottom border (drop left & right) -> borderId 4
$ws2.Range("C1").Borders.Item(7).LineStyle = 0
$ws2.Range("C1").Borders.Item(10).LineStyle = 0

# D1: keep top/right/bottom border (drop left) -> borderId 5
$ws2.Range("D1").Borders.Item(7).LineStyle = 0

# F1: keep only top/bottom border (drop left & right) -> borderId 4
$ws2.Range("F1").Borders.Item(7).LineStyle = 0
$ws2.Range("F1").Borders.Item(10).LineStyle = 0

# G1: keep top/right/bottom border (drop left) -> borderId 5
$ws2.Range("G1").Borders.Item(7).LineStyle = 0

# C2 & F2: rename "fedcore" column headers to "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5: clear the stray empty cell entirely
$ws2.Range("G5").ClearContents()
